# Auto-generated edit script applying Rafflesia_Profits market-price / profit updates
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 94.28570999999999
$ws.Range("I4").Value = 85
$ws.Range("K4").Value = 85
$ws.Range("M4").Value = 29
$ws.Range("H9").Value = 3
$ws.Range("I9").Value = 3
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 166
$ws.Range("H31").Value = 5500
$ws.Range("I31").Value = 5500
$ws.Range("K31").Value = 16500
$ws.Range("M31").Value = -16270
$ws.Range("H38").Value = 6521.4707
$ws.Range("I38").Value = 6521.4707
$ws.Range("K38").Value = 19564.4121
$ws.Range("M38").Value = -19192.4121
$ws.Range("H48").Value = 6339.6665
$ws.Range("J48").Value = 9009.5
$ws.Range("L48").Value = 27028.5
$ws.Range("N48").Value = -27612.5
$ws.Range("H56").Value = 6339.6665
$ws.Range("J56").Value = 9009.5
$ws.Range("L56").Value = 27028.5
$ws.Range("N56").Value = -28096.5
$ws.Range("H82").Value = 22500
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 22500
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H135").Value = 1616.1666
$ws.Range("I135").Value = 1417.8889
$ws.Range("J135").Value = 2211
$ws.Range("K135").Value = 12761.0001
$ws.Range("L135").Value = 19899
$ws.Range("M135").Value = -10226.0001
$ws.Range("N135").Value = -24969
$ws.Range("H138").Value = 1500
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("H141").Value = 1024
$ws.Range("I141").Value = 998.6667
$ws.Range("J141").Value = 1100
$ws.Range("K141").Value = 2996.0001
$ws.Range("L141").Value = 3300
$ws.Range("M141").Value = 2183.9999
$ws.Range("N141").Value = -13660

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 150
$ws.Range("I5").Value = 150
$ws.Range("K5").Value = 150
$ws.Range("M5").Value = -38
$ws.Range("H88").Value = 2237.5
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 2237.5
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 2237.5
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -3049.5
$ws.Range("H91").Value = 2237.5
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 2237.5
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 2237.5
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -5045.5
$ws.Range("H102").Value = 2750
$ws.Range("I102").Value = 2125
$ws.Range("J102").Value = 4000
$ws.Range("K102").Value = 2125
$ws.Range("L102").Value = 4000
$ws.Range("M102").Value = -503
$ws.Range("N102").Value = -7244

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 150
$ws.Range("I4").Value = 150
$ws.Range("K4").Value = 150
$ws.Range("M4").Value = -35
$ws.Range("H36").Value = 7486.6665
$ws.Range("I36").Value = 7486.6665
$ws.Range("K36").Value = 7486.6665
$ws.Range("M36").Value = -6952.6665
$ws.Range("H76").Value = 75662.39999999999
$ws.Range("J76").Value = 75662.39999999999
$ws.Range("L76").Value = 75662.39999999999
$ws.Range("N76").Value = -76292.39999999999
$ws.Range("H79").Value = 75662.39999999999
$ws.Range("J79").Value = 75662.39999999999
$ws.Range("L79").Value = 75662.39999999999
$ws.Range("N79").Value = -77846.39999999999
$ws.Range("H94").Value = 578.1429000000001
$ws.Range("I94").Value = 532.8333
$ws.Range("J94").Value = 850
$ws.Range("K94").Value = 532.8333
$ws.Range("L94").Value = 850
$ws.Range("M94").Value = -81.83330000000001
$ws.Range("N94").Value = -1752
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H107").Value = 3225
$ws.Range("I107").Value = 3213.75
$ws.Range("J107").Value = 3247.5
$ws.Range("K107").Value = 3213.75
$ws.Range("L107").Value = 3247.5
$ws.Range("M107").Value = -1293.75
$ws.Range("N107").Value = -7087.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 50000
$ws.Range("I51").Value = 50000
$ws.Range("K51").Value = 50000
$ws.Range("M51").Value = -49264
$ws.Range("H61").Value = 50000
$ws.Range("I61").Value = 50000
$ws.Range("K61").Value = 50000
$ws.Range("M61").Value = -49652
$ws.Range("H74").Value = 41000
$ws.Range("J74").Value = 41000
$ws.Range("L74").Value = 41000
$ws.Range("N74").Value = -42748
$ws.Range("H77").Value = 41000
$ws.Range("J77").Value = 41000
$ws.Range("L77").Value = 123000
$ws.Range("N77").Value = -131736
$ws.Range("H94").Value = 2208
$ws.Range("I94").Value = 3165.75
$ws.Range("K94").Value = 3165.75
$ws.Range("M94").Value = -2714.75
$ws.Range("H132").Value = 6314.857
$ws.Range("I132").Value = 2644
$ws.Range("J132").Value = 9985.714
$ws.Range("K132").Value = 7932
$ws.Range("L132").Value = 29957.142
$ws.Range("M132").Value = -5402
$ws.Range("N132").Value = -35017.142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 300.5
$ws.Range("I4").Value = 300.5
$ws.Range("K4").Value = 901.5
$ws.Range("M4").Value = -789.5
$ws.Range("H14").Value = 250.5
$ws.Range("I14").Value = 250.5
$ws.Range("K14").Value = 751.5
$ws.Range("M14").Value = -578.5
$ws.Range("H25").Value = 3224.25
$ws.Range("J25").Value = 12000
$ws.Range("L25").Value = 36000
$ws.Range("N25").Value = -36338
$ws.Range("H30").Value = 3224.25
$ws.Range("J30").Value = 12000
$ws.Range("L30").Value = 36000
$ws.Range("N30").Value = -36204
$ws.Range("H51").Value = 250
$ws.Range("I51").Value = 250
$ws.Range("K51").Value = 750
$ws.Range("M51").Value = -290
$ws.Range("H55").Value = 4760
$ws.Range("J55").Value = 5575
$ws.Range("L55").Value = 16725
$ws.Range("N55").Value = -17079
$ws.Range("H138").Value = 1000
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 5002500
$ws.Range("I9").Value = 5002500
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 5002500
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -5002330
$ws.Range("N9").ClearContents()
$ws.Range("H57").Value = 14838.75
$ws.Range("J57").Value = 19777.5
$ws.Range("L57").Value = 19777.5
$ws.Range("N57").Value = -21417.5
$ws.Range("H122").Value = 1319.75
$ws.Range("I122").Value = 1319.75
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3959.25
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1509.25
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 13333.667
$ws.Range("I2").Value = 13333.667
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 13333.667
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -13221.667
$ws.Range("N2").ClearContents()
$ws.Range("H16").Value = 1856.4286
$ws.Range("I16").Value = 1856.4286
$ws.Range("K16").Value = 1856.4286
$ws.Range("M16").Value = -1686.4286
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H136").Value = 23000
$ws.Range("I136").Value = 15666.667
$ws.Range("K136").Value = 47000.001
$ws.Range("M136").Value = -44450.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 50000
$ws.Range("J55").Value = 50000
$ws.Range("L55").Value = 50000
$ws.Range("N55").Value = -50554
$ws.Range("H64").Value = 49995
$ws.Range("I64").Value = 49990
$ws.Range("K64").Value = 49990
$ws.Range("M64").Value = -49742
$ws.Range("H67").Value = 49995
$ws.Range("I67").Value = 49990
$ws.Range("K67").Value = 49990
$ws.Range("M67").Value = -49132
$ws.Range("H93").Value = 50000
$ws.Range("J93").Value = 50000
$ws.Range("L93").Value = 50000
$ws.Range("N93").Value = -54992
$ws.Range("H107").Value = 2425.125
$ws.Range("J107").Value = 2819.4
$ws.Range("L107").Value = 8458.200000000001
$ws.Range("N107").Value = -12298.2
$ws.Range("H126").Value = 5040.75
$ws.Range("I126").Value = 3066.2
$ws.Range("J126").Value = 8331.666999999999
$ws.Range("K126").Value = 9198.599999999999
$ws.Range("L126").Value = 24995.001
$ws.Range("M126").Value = -6728.599999999999
$ws.Range("N126").Value = -29935.001

Write-Host "Applied all Rafflesia_Profits updates."
